$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 133, shifting existing rows
# 133 (and below) down by one.
$ws.Rows(133).Insert()

# Populate the newly inserted row 133 with the new record.
$ws.Range("A133").Value = 7
$ws.Range("B133").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C133").Value = "Ñuble"
$ws.Range("D133").Value = 44628
$ws.Range("E133").Value = 16
$ws.Range("F133").Value = 100112045
$ws.Range("G133").Value = "Zapallo"
$ws.Range("H133").Value = "Camote"
$ws.Range("I133").Value = "1a (cosecha)"
$ws.Range("J133").Value = 600
$ws.Range("K133").Value = 350
$ws.Range("L133").Value = 400
$ws.Range("M133").Value = 375
$ws.Range("N133").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O133").Value = "Región de O'Higgins"
$ws.Range("P133").Value = 375
$ws.Range("Q133").Value = 1
$ws.Range("R133").Value = "Hortaliza"
